$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.374.83"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.98"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.73"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.85"
$ws.Range("E7").Value = "  +5.47%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.642"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.70"
$ws.Range("E10").Value = "  -2.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.07"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.30"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.618.25"
$ws.Range("E15").Value = "  +2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.94"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.874"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.279.25"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.342.96"
$ws.Range("E19").Value = "  +1.10%  "

$ws.Range("E20").Value = "  +3.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.89"
$ws.Range("E22").Value = "  -1.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.67"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +7.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.90"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.65"
$ws.Range("E29").Value = "  -1.24%  "

$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.83"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.01"
$ws.Range("E32").Value = "  +2.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("E33").Value = "  +6.79%  "

$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0812"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.02"
$ws.Range("E36").Value = "  +22.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("E37").Value = "  +2.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").Value = "  +13.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.72"
$ws.Range("E39").Value = "  +2.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0305"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.84"
$ws.Range("E41").Value = "  +13.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.90"
$ws.Range("E43").Value = "  +4.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.214"
$ws.Range("E44").Value = "  +5.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.16"
$ws.Range("E45").Value = "  +6.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.28"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.88"
$ws.Range("E47").Value = "  -4.31%  "

$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.29"
$ws.Range("E51").Value = "  +4.19%  "
